$wb = $excel.ActiveWorkbook

# Reference sheet to borrow the existing bold/centered/bordered header style from
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet after the last existing one so it lands at the end
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "top_3_pays_langues_parless"

# Copy the header formatting (bold font, centered alignment, thin border) from
# an existing header row so we reuse the same style instead of creating a new one
$ws1.Range("A1:C1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

# Header row
$ws.Cells.Item(1, 1).Value = "index"
$ws.Cells.Item(1, 2).Value = "language"
$ws.Cells.Item(1, 3).Value = "population"

# Data rows: top 3 countries by spoken language
$ws.Cells.Item(2, 1).Value = 9
$ws.Cells.Item(2, 2).Value = "Russian"
$ws.Cells.Item(2, 3).Value = 155137605

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Catalan"
$ws.Cells.Item(3, 3).Value = 108331988

$ws.Cells.Item(4, 1).Value = 14
$ws.Cells.Item(4, 2).Value = "German"
$ws.Cells.Item(4, 3).Value = 105240604

$ws.Range("A1").Select() | Out-Null
